# Update "paises" (countries) COVID tracker workbook.
# - Refresh the "last updated" timestamp in A1.
# - Reorder three small-territory rows (Bonaire/Islas Virgenes Britanicas/
#   San Vicente y las Granadinas) and two rows (Montserrat/Islas Malvinas)
#   to match the new source ordering, carrying their labels along.
# - Refresh the numeric COVID statistics (columns B:H) for the updated rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header timestamp
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 14:42"

# ---------------------------------------------------------------------
# Re-label the rows whose country order changed.
#   Before: Islas Virgenes Britanicas / San Vicente y las Granadinas / Bonaire, San Eustaquio y Saba
#   After : Bonaire, San Eustaquio y Saba / Islas Virgenes Britanicas / San Vicente y las Granadinas
# ---------------------------------------------------------------------
$ws.Range("A197").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A198").Value = "Islas Virgenes Britanicas"
$ws.Range("A199").Value = "San Vicente y las Granadinas"

#   Before: Islas Malvinas / Montserrat
#   After : Montserrat / Islas Malvinas
$ws.Range("A214").Value = "Montserrat"
$ws.Range("A215").Value = "Islas Malvinas"

# ---------------------------------------------------------------------
# Helper data: row -> (B, C, D, E, F, G, H)
# Refresh numeric data (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) per updated row.
# ---------------------------------------------------------------------

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Estados Unidos
Set-Row 4 7098674 737 4347494 2545677 0 32 205503

# Kuwait
Set-Row 40 101299 616 92341 8368 0 2 590

# Paises Bajos
Set-Row 41 100597 2357 0 0 0 5 6296

# Suecia
Set-Row 43 89756 0 0 0 0 5 5876

# Bielorrusia
Set-Row 49 76357 253 73564 1997 0 5 796

# Suiza
Set-Row 61 51101 437 42100 6941 0 6 2060

# Estado de Palestina
Set-Row 71 37083 503 26288 10523 0 3 272

# Dinamarca
Set-Row 80 24357 558 18359 5355 0 2 643

# Consejo Danes para los Refugiados
Set-Row 97 10537 14 10041 225 0 0 271

# Uganda
Set-Row 113 6879 167 2961 3849 0 5 69

# Sri Lanka
Set-Row 142 3315 2 3129 173 0 0 13

# Gibraltar
Set-Row 182 357 2 324 33 0 0 0

# Liechtenstein
Set-Row 195 116 1 110 5 0 0 1

# Row 197 -> now "Bonaire, San Eustaquio y Saba"
Set-Row 197 69 15 21 47 0 0 1

# Row 198 -> now "Islas Virgenes Britanicas"
Set-Row 198 69 0 48 20 0 0 1

# Row 199 -> now "San Vicente y las Granadinas"
Set-Row 199 64 0 64 0 0 0 0

# Row 214 -> now "Montserrat"
Set-Row 214 13 0 12 0 0 0 1

# Row 215 -> now "Islas Malvinas"
Set-Row 215 13 0 13 0 0 0 0
